$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 3.239412069320679
$ws.Range("B1").Value = 2.448537826538086
$ws.Range("C1").Value = 2.00375771522522
$ws.Range("D1").Value = 2.029806852340698
$ws.Range("E1").Value = 2.294777393341064
